# Agregados IDs y creadas relaciones en la tabla
#
# - Sheet "obras": adds a leading "ID" column (1,2,3) and collapses the three
#   "referente" columns (Titulo referente / Fecha referente / Periodico /
#   archivo referente) into a single "Referentes" column that stores the
#   numeric id(s) of the matching row(s) in the "refentes" sheet.
# - Sheet "refentes": adds a leading "ID" column (1,2,3,4) and a new 4th
#   referente row ("Exmilitar Mata a la Esposa de su Amigo y se Suicida",
#   with the spelling of "Esposa" corrected) that used to be crammed into
#   the "obras" sheet's row 4 referente columns.

$wb = $excel.ActiveWorkbook

$obras = $wb.Worksheets.Item("obras")
$refentes = $wb.Worksheets.Item("refentes")

# ---------------------------------------------------------------------
# Sheet "refentes": insert the ID column, number the existing rows, and
# append the new 4th referente (moved out of "obras").
# ---------------------------------------------------------------------
$refentes.Columns.Item(1).Insert()

$refentes.Cells.Item(1, 1).Value = "ID"
$refentes.Cells.Item(2, 1).Value = 1
$refentes.Cells.Item(3, 1).Value = 2
$refentes.Cells.Item(4, 1).Value = 3
$refentes.Cells.Item(5, 1).Value = 4

$refentes.Cells.Item(5, 2).Value = "Exmilitar Mata a la Esposa de su Amigo y se Suicida"
$refentes.Cells.Item(5, 5).Value = "exmilitar-mata-esposa.jpg"

# ---------------------------------------------------------------------
# Sheet "obras": insert the ID column, then replace the old referente
# columns (now F:I) with a single "Referentes" column holding the row
# id(s) from "refentes".
# ---------------------------------------------------------------------
$obras.Columns.Item(1).Insert()
$obras.Range("G1:J1").EntireColumn.Delete()

$obras.Cells.Item(1, 1).Value = "ID"
$obras.Cells.Item(2, 1).Value = 1
$obras.Cells.Item(3, 1).Value = 2
$obras.Cells.Item(4, 1).Value = 3

$obras.Cells.Item(1, 7).Value = "Referentes"
$obras.Cells.Item(2, 7).Value = 1
$obras.Cells.Item(3, 7).Value = "2,3"
$obras.Cells.Item(4, 7).Value = 4

# ---------------------------------------------------------------------
# View state: "obras" becomes the active/selected sheet (was "refentes"),
# with its selection on the new Referentes column; "refentes" keeps a
# stale-looking selection further down the sheet.
# ---------------------------------------------------------------------
$refentes.Activate()
$refentes.Range("C10").Select()

$obras.Activate()
$obras.Range("G3").Select()
